$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 203, shifting existing rows 203:234 down to 204:235
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new weekly record.
# Columns A, B, C, E, F, G, H, I, O, R are constant across this data block,
# so re-use the same literal values found in the surrounding rows.
$ws.Range("A203").Value = 10
$ws.Range("B203").Value = "Vega Modelo de Temuco"
$ws.Range("C203").Value = "La Araucanía"
$ws.Range("D203").Value = 45180
$ws.Range("E203").Value = 9
$ws.Range("F203").Value = 100114002
$ws.Range("G203").Value = "Camote"
$ws.Range("H203").Value = "Sin especificar"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 80
$ws.Range("K203").Value = 24000
$ws.Range("L203").Value = 24000
$ws.Range("M203").Value = 24000
$ws.Range("N203").Value = "$/caja 18 kilos"
$ws.Range("O203").Value = "Perú"
$ws.Range("P203").Value = 1333
$ws.Range("Q203").Value = 18
$ws.Range("R203").Value = "Hortaliza"

# Give the new row's date cell the same date number-format style as the rest
# of column D in this block.
$ws.Range("D203").NumberFormat = $ws.Range("D204").NumberFormat()
